# Applies the "4wk low sales check" update to the forecast workbook:
#  - Updates Seasonality Index (column L) values on the "Forecast Comparison" sheet
#  - Marks the last 4 weeks (rows 14-17) as a low-sales forecast:
#      MyForecast (D) -> 0, Inventory Coverage (H) -> blank,
#      Stockout Risk (I) -> "Low", Reorder Urgency (J) -> "Normal"
#  - Refreshes the dependent totals on the "Summary" sheet

$wb = $excel.ActiveWorkbook

$fc = $wb.Worksheets.Item("Forecast Comparison")

# --- Seasonality Index (column L) updates for every week row ---
$seasonality = @{
    2  = 1.11
    3  = 1.14
    4  = 1.08
    5  = 1.07
    6  = 1.17
    7  = 1.01
    8  = 1.16
    9  = 1.02
    10 = 0.9
    11 = 0.88
    12 = 0.97
    13 = 1.03
    14 = 0.8
    15 = 0.82
    16 = 1.14
    17 = 1.07
}

foreach ($row in $seasonality.Keys) {
    $fc.Range("L$row").Value = $seasonality[$row]
}

# --- Rows 14-17 (the most recent 4 weeks): mark as low-sales weeks ---
foreach ($row in 14..17) {
    $fc.Range("D$row").Value = 0
    $fc.Range("H$row").Value = ""
    $fc.Range("I$row").Value = "Low"
    $fc.Range("J$row").Value = "Normal"
}

# --- Update the Summary sheet with the refreshed totals ---
# (Summary!B column stores these totals as text, so a leading apostrophe is
#  used to force a text value instead of Excel auto-converting to a number;
#  the style is then reset so no extra number formatting is left behind.)
$summary = $wb.Worksheets.Item("Summary")

function Set-TextValue($range, $text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

Set-TextValue $summary.Range("B9")  "12"
Set-TextValue $summary.Range("B10") "8"
Set-TextValue $summary.Range("B12") "1"
Set-TextValue $summary.Range("B14") "0"
